# This workbook contains a small daily price table for "Corazón de apio".
# The edit cyclically shifts the data rows (2-5) down by one position:
#   old row 5 -> new row 2
#   old row 2 -> new row 3
#   old row 3 -> new row 4
#   old row 4 -> new row 5
# Columns A, B, C, E, F, G, H, O, R are identical across all rows, so only
# columns D, I, J, K, L, M, N, P, Q need to be rewritten per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values for the columns that actually change, for rows 2-5.
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

$rows = @{}
foreach ($r in 2..5) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $rows[$r] = $rowData
}

# New order after shifting down by one (with wraparound): row 2 gets old row 5's
# values, row 3 gets old row 2's, row 4 gets old row 3's, row 5 gets old row 4's.
$newOrder = @{
    2 = $rows[5]
    3 = $rows[2]
    4 = $rows[3]
    5 = $rows[4]
}

foreach ($r in 2..5) {
    $rowData = $newOrder[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value2 = $rowData[$col]
    }
}
